$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Update the "time_taken" timestamps in column F of the "data" sheet.
# ---------------------------------------------------------------------------
$dataWs.Range("F2").Value  = "2021-10-05 14:20:10.947110"
$dataWs.Range("F3").Value  = "2021-10-05 14:20:10.947118"
$dataWs.Range("F4").Value  = "2021-10-05 14:20:10.947121"
$dataWs.Range("F5").Value  = "2021-10-05 14:20:10.947124"
$dataWs.Range("F6").Value  = "2021-10-05 14:20:10.947127"
$dataWs.Range("F7").Value  = "2021-10-05 14:20:10.947129"
$dataWs.Range("F8").Value  = "2021-10-05 14:20:10.947132"
$dataWs.Range("F9").Value  = "2021-10-05 14:20:10.947134"
$dataWs.Range("F10").Value = "2021-10-05 14:20:10.947137"
$dataWs.Range("F11").Value = "2021-10-05 14:20:10.947139"
$dataWs.Range("F12").Value = "2021-10-05 14:20:10.947142"
$dataWs.Range("F13").Value = "2021-10-05 14:20:10.947144"
$dataWs.Range("F14").Value = "2021-10-05 14:20:10.947147"
$dataWs.Range("F15").Value = "2021-10-05 14:20:10.947149"
$dataWs.Range("F16").Value = "2021-10-05 14:20:10.947152"
$dataWs.Range("F17").Value = "2021-10-05 14:20:10.947154"

# ---------------------------------------------------------------------------
# 2. Add the new "metadata" worksheet after "data".
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $dataWs)
$ws2.Name = "metadata"

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Familial cerebral small vessel disease"
$ws2.Range("C2").Value = 50
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "1.10"
$ws2.Range("E2").Value = "2020-10-06T15:57:08.085104Z"
$ws2.Range("F2").Value = "2021-10-05 14:20:10.943608"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/50/?format=json"

# Match the header / index-column formatting used on the "data" sheet
# (bold, centered, top-aligned, thin border) by copying the existing style.
$dataWs.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$dataWs.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$dataWs.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$excel.CutCopyMode = 0
